# Auto-generated edit script to apply the scraped-data refresh diff
$wb = $excel.ActiveWorkbook

### Sheet 1 ###
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = "Última actualización: 20:53:28"
$ws.Range("A3").Value = "Total filas: 371"
$ws.Range("C38").Value = "11_ETCHEVERRY"
$ws.Range("C39").Value = "15_ABASTO"
$ws.Range("A62").Value = "08:38:24"
$ws.Range("C62").Value = "27_EL RETIRO"
$ws.Range("D62").Value = 39
$ws.Range("A63").Value = "07:49:32"
$ws.Range("C63").Value = "14_ABASTO"
$ws.Range("D63").Value = 88
$ws.Range("A178").Value = "12:33:02"
$ws.Range("C178").Value = "27_EL RETIRO"
$ws.Range("D178").Value = 104
$ws.Range("A179").Value = "12:53:26"
$ws.Range("C179").Value = "11_ETCHEVERRY"
$ws.Range("D179").Value = 84
$ws.Range("A190").Value = "13:55:43"
$ws.Range("C190").Value = "16_SANTA ANA"
$ws.Range("D190").Value = 56
$ws.Range("A191").Value = "13:41:21"
$ws.Range("C191").Value = "23_HERNANDEZ"
$ws.Range("D191").Value = 70
$ws.Range("A235").Value = "16:44:58"
$ws.Range("C235").Value = "16_P MOR-SANTA ANA"
$ws.Range("D235").Value = 22
$ws.Range("A236").Value = "16:28:21"
$ws.Range("C236").Value = "23_HERNANDEZ"
$ws.Range("D236").Value = 38
$ws.Range("A247").Value = "17:35:41"
$ws.Range("C247").Value = "215B_EL PATO"
$ws.Range("D247").Value = 2
$ws.Range("A248").Value = "16:12:06"
$ws.Range("C248").Value = "27_EL RETIRO"
$ws.Range("D248").Value = 85
$ws.Range("A258").Value = "16:51:51"
$ws.Range("C258").Value = "10_OLMOS"
$ws.Range("D258").Value = 62
$ws.Range("A259").Value = "16:37:37"
$ws.Range("C259").Value = "23_HERNANDEZ"
$ws.Range("D259").Value = 76
$ws.Range("C269").Value = "16_P MOR-SANTA ANA"
$ws.Range("C270").Value = "15_ABASTO"
$ws.Range("A294").Value = "17:35:41"
$ws.Range("C294").Value = "215_EL PELIGRO"
$ws.Range("D294").Value = 88
$ws.Range("A295").Value = "17:55:25"
$ws.Range("C295").Value = "14_ABASTO"
$ws.Range("D295").Value = 68
$ws.Range("A351").Value = "20:32:11"
$ws.Range("C351").Value = "14_ABASTO"
$ws.Range("D351").Value = 37
$ws.Range("A352").Value = "19:47:50"
$ws.Range("C352").Value = "15_ABASTO"
$ws.Range("D352").Value = 82
$ws.Range("A371").Value = "20:53:28"
$ws.Range("B371").Value = "22:33"
$ws.Range("D371").Value = 100
$ws.Range("C372").Value = "215C_EL PATO"
$ws.Range("B373").Value = "22:34"
$ws.Range("C373").Value = "16_SANTA ANA"
$ws.Range("D373").Value = 108
$ws.Range("A374").Value = "20:53:28"
$ws.Range("B374").Value = "22:35"
$ws.Range("C374").Value = "23_HERNANDEZ"
$ws.Range("D374").Value = 102
$ws.Range("E374").Value = "LP1912"
$ws.Range("A375").Value = "20:53:28"
$ws.Range("B375").Value = "22:43"
$ws.Range("C375").Value = "11X44_ETCHEVERRY"
$ws.Range("D375").Value = 110
$ws.Range("E375").Value = "LP1912"
$ws.Range("A376").Value = "20:46:15"
$ws.Range("B376").Value = "22:44"
$ws.Range("C376").Value = "11X44_ETCHEVERRY"
$ws.Range("D376").Value = 118
$ws.Range("E376").Value = "LP1912"

### Sheet 2 ###
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = "Última actualización: 20:53:28"
$ws.Range("A3").Value = "Total filas: 56"
$ws.Range("A60").Value = "20:53:28"
$ws.Range("B60").Value = "22:33"
$ws.Range("D60").Value = 100
$ws.Range("A61").Value = "20:46:15"
$ws.Range("B61").Value = "22:34"
$ws.Range("C61").Value = "215C_EL PATO"
$ws.Range("D61").Value = 108
$ws.Range("E61").Value = "LP1912"

### Sheet 3 ###
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = "Última actualización: 20:53:28"
$ws.Range("A3").Value = "Total filas: 52"
$ws.Range("A53").Value = "20:53:28"
$ws.Range("B53").Value = "20:53"
$ws.Range("D53").Value = 0
$ws.Range("A54").Value = "19:47:50"
$ws.Range("B54").Value = "21:27"
$ws.Range("D54").Value = 100
$ws.Range("A55").Value = "19:54:57"
$ws.Range("B55").Value = "21:29"
$ws.Range("D55").Value = 95
$ws.Range("A56").Value = "19:35:34"
$ws.Range("B56").Value = "21:30"
$ws.Range("C56").Value = "215C_LA PLATA"
$ws.Range("D56").Value = 115
$ws.Range("E56").Value = "L6203"
$ws.Range("A57").Value = "20:32:11"
$ws.Range("B57").Value = "22:20"
$ws.Range("C57").Value = "215B_LP-P MOR-40 Y 115"
$ws.Range("D57").Value = 108
$ws.Range("E57").Value = "L6173"
